# Update Name of Algo
# Apply targeted numeric corrections to the KNN-imputed result data on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6.228
$ws.Range("A9").Value = -20.912
$ws.Range("A18").Value = -21.825
$ws.Range("A20").Value = -21.738
$ws.Range("C21").Value = -12.688

$wb.Save()
